$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) stores values as literal text (e.g. "27.062.83",
# "1.878.44", "0.000008201") rather than numbers, since several values use
# a thousands-dot format or would otherwise lose significant trailing
# zeros/precision if Excel auto-converted them to numeric cells. Force
# each updated Price cell to Text format before writing its new value so
# the literal string is preserved exactly as scraped from the source.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.062.83"
$ws.Range("E2").Value = "  +5.41%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.878.44"
$ws.Range("E3").Value = "  +4.01%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "280.89"
$ws.Range("E5").Value = "  +2.24%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9997"
$ws.Range("E6").Value = "  -0.01%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5282"
$ws.Range("E7").Value = "  +4.65%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3522"
$ws.Range("E8").Value = "  +0.36%  "
$ws.Range("E9").Value = "  +2.47%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07037"
$ws.Range("E10").Value = "  +6.38%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.34"
$ws.Range("E11").Value = "  +2.39%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.8153"
$ws.Range("E12").Value = "  -2.03%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07805"
$ws.Range("E13").Value = "  +0.20%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.882.50"
$ws.Range("E14").Value = "  +4.21%  "
$ws.Range("E15").Value = "  +2.98%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "90.48"
$ws.Range("E16").Value = "  +3.59%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.000"
$ws.Range("E17").Value = "  +0.06%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.58"
$ws.Range("E18").Value = "  +5.08%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000008201"
$ws.Range("E19").Value = "  +2.96%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9999"
$ws.Range("E20").Value = "  +0.02%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "27.090.23"
$ws.Range("E21").Value = "  +5.28%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.116.06"
$ws.Range("E22").Value = "  +4.12%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.766"
$ws.Range("E23").Value = "  +1.17%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.17"
$ws.Range("E24").Value = "  +2.47%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "6.225"
$ws.Range("E25").Value = "  +3.00%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.387"
$ws.Range("E26").Value = "  +12.91%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "146.49"
$ws.Range("E27").Value = "  +3.29%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.57"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.674"
$ws.Range("E29").Value = "  +1.45%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "112.58"
$ws.Range("E30").Value = "  +3.94%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.394"
$ws.Range("E31").Value = "  +1.77%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.375"
$ws.Range("E32").Value = "  +4.59%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.08907"
$ws.Range("E33").Value = "  +1.56%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04905"
$ws.Range("E34").Value = "  +2.52%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.175"
$ws.Range("E35").Value = "  +3.82%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7457"
$ws.Range("E36").Value = "  +3.25%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.894"
$ws.Range("E37").Value = "  +0.88%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.306"
$ws.Range("E38").Value = "  +9.10%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.412"
$ws.Range("E39").Value = "  +6.68%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5315"
$ws.Range("E40").Value = "  +2.86%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.01884"
$ws.Range("E41").Value = "  +1.80%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9777"
$ws.Range("E42").Value = "  +3.43%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "117.15"
$ws.Range("E43").Value = "  +4.13%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "6.323"
$ws.Range("E44").Value = "  +2.98%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.193"
$ws.Range("E45").Value = "  +2.63%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.9988"
$ws.Range("E46").Value = "  -0.05%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4602"
$ws.Range("E47").Value = "  +1.27%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1370"
$ws.Range("E48").Value = "  -0.27%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.504"
$ws.Range("E49").Value = "  +2.21%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "36.72"
$ws.Range("E50").Value = "  +2.11%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.529"
$ws.Range("E51").Value = "  +2.89%  "
